# Update gh-pages to output generated at 456a3b4
# Applies updated "想去人数" (F) and "最低票价" (G) figures to the
# "展览" and "全部类型" worksheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" -----------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("F2").Value = 5712
$ws1.Range("F3").Value = 85
$ws1.Range("G3").Value = 50
$ws1.Range("F6").Value = 159
$ws1.Range("F7").Value = 2616
$ws1.Range("F8").Value = 87
$ws1.Range("F9").Value = 188
$ws1.Range("F11").Value = 98
$ws1.Range("F12").Value = 39
$ws1.Range("F13").Value = 2456
$ws1.Range("F14").Value = 501

# --- Sheet "全部类型" ---------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F2").Value = 5712
$ws4.Range("F3").Value = 85
$ws4.Range("G3").Value = 50
$ws4.Range("F7").Value = 159
$ws4.Range("F8").Value = 2616
$ws4.Range("F9").Value = 87
$ws4.Range("F10").Value = 188
$ws4.Range("F13").Value = 98
$ws4.Range("F14").Value = 39
$ws4.Range("F15").Value = 2456
$ws4.Range("F16").Value = 501
